# Rename the two Pearson logo pictures (in the "first" and "default"
# footers) from "image2.png" to "image1.png", and the BTEC logo picture
# (in the "first" header) from "image1.jpg" to "image2.jpg".
#
# Word's InlineShape object has no settable Name property directly, so
# each picture is momentarily converted to a floating Shape (which does
# expose .Name), renamed, then converted back to an inline picture so the
# layout/anchoring is unchanged.

$d = $word.ActiveDocument

foreach ($story in $d.StoryRanges) {
    $shapeCount = $story.InlineShapes.Count
    for ($i = 1; $i -le $shapeCount; $i++) {
        $shp = $story.InlineShapes.Item($i).ConvertToShape()
        $currentName = $shp.Name

        if ($currentName -eq "image2.png") {
            $shp.Name = "image1.png"
        } elseif ($currentName -eq "image1.jpg") {
            $shp.Name = "image2.jpg"
        }

        $shp.ConvertToInlineShape() | Out-Null
    }
}
